{"js": "const body = context.document.body;\n\n// The paragraph currently reads \"Version 2.\" (with \"Version\" split across\n// two runs - \"Versi\" + \"on\" - left over from an earlier edit). Normalize\n// that split into a single \"Version\" run, matching how Word coalesces\n// same-formatted runs when the paragraph is next touched.\nconst versionResults = body.search(\"Version\", { matchCase: true, matchWholeWord: false });\nversionResults.load(\"items\");\nawait context.sync();\nif (versionResults.items.length > 0) {\n  versionResults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Bump the version number from 2 to 1, folding the final period into the\n// same run (the trailing \".\" used to live in its own run after the\n// \"_GoBack\" bookmark).\nconst numberResults = body.search(\" 2\", { matchCase: true, matchWholeWord: false });\nnumberResults.load(\"items\");\nawait context.sync();\nif (numberResults.items.length > 0) {\n  numberResults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Remove the now-redundant standalone trailing period run.\nconst periodResults = body.search(\".\", { matchCase: true, matchWholeWord: false });\nperiodResults.load(\"items\");\nawait context.sync();\nif (periodResults.items.length > 0) {\n  periodResults.items[periodResults.items.length - 1].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The paragraph currently reads \"Version 2.\" (with \"Version\" split across\n# two runs - \"Versi\" + \"on\" - left over from an earlier edit). Normalize\n# that split into a single \"Version\" run, matching how Word coalesces\n# same-formatted runs when the paragraph is next touched. (Setting the\n# range's Text to its own value is a no-op, so nudge it through a\n# temporary value first to force the rewrite.)\n$versionRange = $d.Content\n$null = $versionRange.Find.Execute(\"Version\")\n$versionRange.Text = \"Versionx\"\n$versionRange2 = $d.Content\n$null = $versionRange2.Find.Execute(\"Versionx\")\n$versionRange2.Text = \"Version\"\n\n# Bump the version number from 2 to 1, folding the final period into the\n# same run (the trailing \".\" used to live in its own run after the\n# \"_GoBack\" bookmark).\n$numberRange = $d.Content\n$null = $numberRange.Find.Execute(\" 2\")\n$numberRange.Text = \" 1.\"\n\n# Remove the now-redundant standalone trailing period run.\n$charCount = $d.Characters.Count\n$trailingDot = $d.Range($charCount - 2, $charCount - 1)\nif ($trailingDot.Text -eq \".\") {\n    $trailingDot.Text = \"\"\n}\n"}
